# Update the Questa vplan workbook: append a new "test status" section
# (row 26) to the Testplan sheet, matching the rows above it in style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testplan")
$ws.Activate()

# Row 26 reuses the same cell formatting as row 25 (the last data row) -
# copy its formats down first so the new row's style indices line up with
# the existing ones instead of minting new cellXfs entries.
$ws.Range("B25:H25").Copy() | Out-Null
$ws.Range("B26").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Section number column (B) holds text like "1.1", "6", "7" elsewhere in
# the sheet, so force "8" to be stored as text too (leading apostrophe is
# the classic "type this as text" marker and keeps the copied style).
$ws.Range("B26").Value = "'8"
$ws.Range("C26").Value = "test status"
$ws.Range("D26").Value = "All tests status and pass rate recorded"
$ws.Range("E26").Value = "kei_i2c_reg_access_test_*`nkei_i2c_reg_bit_bash_test_*`nkei_i2c_reg_hw_reset_test_*       `nkei_i2c_quick_reg_access_test_*`nkei_i2c_master_directed_write_packet_test_*`nkei_i2c_master_directed_read_packet_test_*`nkei_i2c_master_directed_interrupt_test_*`nkei_i2c_master_address_cg_test_*`nkei_i2c_master_ss_cnt_test_*`nkei_i2c_master_fs_cnt_test_*`nkei_i2c_master_hs_cnt_test_*`nkei_i2c_master_sda_control_cg_test_*`nkei_i2c_master_timeout_cg_test_*`nkei_i2c_master_enabled_cg_test_*`nkei_i2c_master_stop_det_intr_test_*`nkei_i2c_master_tx_abrt_intr_test_*`nkei_i2c_master_rx_full_intr_test_*`nkei_i2c_master_rx_over_intr_test_*`n"
$ws.Range("F26").Value = "Test"
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 100

# Match the source row height for the new, much taller wrapped-text row.
$ws.Rows.Item(26).RowHeight = 268.75

# Move the active selection down to the new row, like the saved workbook.
$ws.Range("E26").Select() | Out-Null
